$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 183
$lastExistingRow = 182
$numNewRows = 7
$lastNewRow = $firstNewRow + $numNewRows - 1

$url = "https://github.com/ersilia-os/ersilia"
$repoName = "ersilia"
$repoAuthor = "ersilia-os"
$startDate = "07/04/2020"
# OSE BCE PDE SV OS SD RS TFS UI TC
$vals = @("0","0","1","1","1","0","0","0","0","1")

# Match the bold/border/centered look already used by the id column (column A)
# on the preceding data rows by copying that cell's format onto the new ids.
for ($i = 0; $i -lt $numNewRows; $i++) {
    $row = $firstNewRow + $i
    $ws.Cells.Item($lastExistingRow, 1).Copy($ws.Cells.Item($row, 1))
}

# Text columns (B:O) must stay plain text, not get auto-coerced into numbers/dates,
# so format as text before writing, then drop the temporary text format afterwards
# so no stray formatting is left behind (matches the plain/general cells elsewhere).
$textRange = $ws.Range("B" + $firstNewRow + ":O" + $lastNewRow)
$textRange.NumberFormat = "@"

for ($i = 0; $i -lt $numNewRows; $i++) {
    $row = $firstNewRow + $i
    $id = $lastExistingRow + $i

    $ws.Cells.Item($row, 1).Value = $id
    $ws.Cells.Item($row, 2).Value = $url
    $ws.Cells.Item($row, 3).Value = $repoName
    $ws.Cells.Item($row, 4).Value = $repoAuthor
    $ws.Cells.Item($row, 5).Value = $startDate

    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, 6 + $c).Value = $vals[$c]
    }
}

$textRange.ClearFormats()
